$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '26.416.62'
Set-TextValue 'E2' '  -7.42%  '
Set-TextValue 'D3' '1.687.97'
Set-TextValue 'E3' '  -5.88%  '
Set-TextValue 'D4' '1.004'
Set-TextValue 'E4' '  +0.11%  '
Set-TextValue 'D5' '219.39'
Set-TextValue 'E5' '  -5.19%  '
Set-TextValue 'D6' '0.5093'
Set-TextValue 'E6' '  -13.68%  '
Set-TextValue 'D7' '1.004'
Set-TextValue 'E7' '  +0.05%  '
Set-TextValue 'D8' '0.2681'
Set-TextValue 'E8' '  -3.24%  '
Set-TextValue 'D9' '22.09'
Set-TextValue 'E9' '  -5.58%  '
Set-TextValue 'D10' '0.06316'
Set-TextValue 'E10' '  -6.48%  '
Set-TextValue 'D11' '0.07394'
Set-TextValue 'E11' '  -2.19%  '
Set-TextValue 'D12' '1.686.66'
Set-TextValue 'E12' '  -6.01%  '
Set-TextValue 'D13' '4.538'
Set-TextValue 'E13' '  -5.23%  '
Set-TextValue 'D14' '0.5796'
Set-TextValue 'E14' '  -5.21%  '
Set-TextValue 'D15' '1.914.97'
Set-TextValue 'E15' '  -5.95%  '
Set-TextValue 'D16' '0.000008639'
Set-TextValue 'E16' '  -2.64%  '
Set-TextValue 'D17' '65.20'
Set-TextValue 'E17' '  -13.85%  '
Set-TextValue 'D18' '26.463.72'
Set-TextValue 'E18' '  -7.25%  '
Set-TextValue 'D19' '5.003'
Set-TextValue 'E19' '  -7.65%  '
Set-TextValue 'E20' '  +0.18%  '
Set-TextValue 'D21' '10.90'
Set-TextValue 'E21' '  -4.88%  '
Set-TextValue 'D22' '186.53'
Set-TextValue 'E22' '  -10.57%  '
Set-TextValue 'D23' '6.267'
Set-TextValue 'E23' '  -8.10%  '
Set-TextValue 'D24' '1.004'
Set-TextValue 'E24' '  +0.05%  '
Set-TextValue 'D25' '144.74'
Set-TextValue 'E25' '  -5.16%  '
Set-TextValue 'D26' '7.512'
Set-TextValue 'E26' '  -6.05%  '
Set-TextValue 'D27' '0.1174'
Set-TextValue 'E27' '  -7.11%  '
Set-TextValue 'D28' '15.87'
Set-TextValue 'E28' '  -3.38%  '
Set-TextValue 'D29' '1.351'
Set-TextValue 'E29' '  -4.65%  '
Set-TextValue 'D30' '0.05759'
Set-TextValue 'E30' '  -5.83%  '
Set-TextValue 'D31' '1.335'
Set-TextValue 'E31' '  -6.03%  '
Set-TextValue 'D32' '3.527'
Set-TextValue 'E32' '  -6.81%  '
Set-TextValue 'D33' '3.526'
Set-TextValue 'E33' '  -6.23%  '
Set-TextValue 'D34' '1.662'
Set-TextValue 'E34' '  -3.75%  '
Set-TextValue 'D35' '1.016'
Set-TextValue 'E35' '  -3.15%  '
Set-TextValue 'D36' '0.5970'
Set-TextValue 'E36' '  -6.98%  '
Set-TextValue 'D37' '2.354'
Set-TextValue 'E37' '  -5.94%  '
Set-TextValue 'D38' '2.678'
Set-TextValue 'E38' '  -0.93%  '
Set-TextValue 'D39' '1.103.69'
Set-TextValue 'E39' '  -3.99%  '
Set-TextValue 'D40' '0.01616'
Set-TextValue 'E40' '  -4.28%  '
Set-TextValue 'D41' '5.896'
Set-TextValue 'E41' '  -6.66%  '
Set-TextValue 'D42' '0.8623'
Set-TextValue 'E42' '  -1.29%  '
Set-TextValue 'D43' '1.003'
Set-TextValue 'E43' '  -0.01%  '
Set-TextValue 'D44' '99.90'
Set-TextValue 'E44' '  -0.56%  '
Set-TextValue 'D45' '1.840.29'
Set-TextValue 'E45' '  -5.51%  '
Set-TextValue 'D46' '0.00000000112'
Set-TextValue 'E46' '  +2.40%  '
Set-TextValue 'D47' '56.51'
Set-TextValue 'E47' '  -6.13%  '
Set-TextValue 'D48' '1.004'
Set-TextValue 'E48' '  +0.46%  '
Set-TextValue 'D49' '8.053'
Set-TextValue 'E49' '  -3.18%  '
Set-TextValue 'B50' 'Cronos'
Set-TextValue 'C50' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D50' '0.05219'
Set-TextValue 'E50' '  -4.24%  '
Set-TextValue 'B51' 'Mantle'
Set-TextValue 'C51' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D51' '0.4309'
Set-TextValue 'E51' '  -3.68%  '

Write-Output "Applied 103 cell updates"
